# Update crypto price (Price) and volume-change (Volume(1h)) columns
# with the latest scraped values. D-column cells are forced to Text
# format first so values like '27.252.13' / '0.9987' stay text (as
# in the source data) instead of being auto-converted to numbers/dates.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '27.252.13'
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.908.04'
$ws.Range('E3').Value = '  +2.11%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.9987'
$ws.Range('E4').Value = '  -0.22%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '307.91'
$ws.Range('E5').Value = '  +0.95%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.9994'
$ws.Range('E6').Value = '  -0.20%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.5242'
$ws.Range('E7').Value = '  +3.03%  '
$ws.Range('E8').Value = '  +3.67%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.07306'
$ws.Range('E9').Value = '  +1.73%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '21.33'
$ws.Range('E10').Value = '  +3.11%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.9009'
$ws.Range('E11').Value = '  +0.98%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.07672'
$ws.Range('E12').Value = '  +2.02%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.903.74'
$ws.Range('E13').Value = '  +1.82%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '95.25'
$ws.Range('E14').Value = '  +0.43%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '5.258'
$ws.Range('E15').Value = '  +0.72%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.9989'
$ws.Range('E16').Value = '  -0.27%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.000008608'
$ws.Range('E17').Value = '  +1.33%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '14.53'
$ws.Range('E18').Value = '  +2.45%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.9987'
$ws.Range('E19').Value = '  -0.08%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '27.299.71'
$ws.Range('E20').Value = '  +1.57%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '5.087'
$ws.Range('E21').Value = '  +1.45%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '2.149.72'
$ws.Range('E22').Value = '  +1.52%  '
$ws.Range('E23').Value = '  +2.72%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '6.445'
$ws.Range('E24').Value = '  +1.14%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.317'
$ws.Range('E25').Value = '  +10.95%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '145.94'
$ws.Range('E26').Value = '  -1.47%  '
$ws.Range('E27').Value = '  +1.75%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.730'
$ws.Range('E28').Value = '  -3.00%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '114.87'
$ws.Range('E29').Value = '  +1.34%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '4.959'
$ws.Range('E30').Value = '  +5.26%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '4.819'
$ws.Range('E31').Value = '  +2.44%  '
$ws.Range('E32').Value = '  +0.86%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.05082'
$ws.Range('E33').Value = '  +0.15%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.248'
$ws.Range('E34').Value = '  +8.04%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.7925'
$ws.Range('E35').Value = '  +6.10%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.992'
$ws.Range('E36').Value = '  +0.33%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '3.300'
$ws.Range('E37').Value = '  +2.27%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.622'
$ws.Range('E38').Value = '  +3.84%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.5681'
$ws.Range('E39').Value = '  +1.14%  '
$ws.Range('E40').Value = '  -0.12%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.074'
$ws.Range('E41').Value = '  +0.02%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '9.026'
$ws.Range('E42').Value = '  +5.43%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '6.639'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '119.03'
$ws.Range('E44').Value = '  +3.25%  '
$ws.Range('E45').Value = '  +3.26%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.4859'
$ws.Range('E46').Value = '  +2.68%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '10.21'
$ws.Range('E47').Value = '  +0.85%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.9991'
$ws.Range('E48').Value = '  -0.27%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.607'
$ws.Range('E49').Value = '  +2.55%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '37.45'
$ws.Range('E50').Value = '  +1.55%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '64.36'
$ws.Range('E51').Value = '  +2.09%  '
